# Append 75 new rows (53-127) of exam room data to the "PT" sheet,
# following the existing pattern: column A = "B0NN", column B = "L0NN",
# column C = "K" + N (no leading zeros), for N = 52..126.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($n = 52; $n -le 126; $n++) {
    $row = $n + 1
    $padded = $n.ToString().PadLeft(3, '0')
    $ws.Cells.Item($row, 1).Value = "B" + $padded
    $ws.Cells.Item($row, 2).Value = "L" + $padded
    $ws.Cells.Item($row, 3).Value = "K" + $n
}

# Match the saved selection state: activeCell A52, selected range A52:C127.
[void]$ws.Range("A52:C127").Select()
